$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74

# Date / text-ish columns: force literal text (no date/number auto-parsing),
# then strip the number-format override so the cell keeps the sheet's
# default (unstyled) look, matching the other data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-24"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "22:26:49"

$ws.Cells.Item($row, 3).Value = "Saturday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "25"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 122597
$ws.Cells.Item($row, 6).Value = 134126
$ws.Cells.Item($row, 7).Value = 163050
$ws.Cells.Item($row, 8).Value = 133459
$ws.Cells.Item($row, 9).Value = 177613
$ws.Cells.Item($row, 10).Value = 115843
$ws.Cells.Item($row, 11).Value = 203167
$ws.Cells.Item($row, 12).Value = 226178
$ws.Cells.Item($row, 13).Value = 175684
$ws.Cells.Item($row, 14).Value = 104251
$ws.Cells.Item($row, 15).Value = 39582
$ws.Cells.Item($row, 16).Value = 33818
$ws.Cells.Item($row, 17).Value = 52004
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36186
$ws.Cells.Item($row, 20).Value = -1
